# Reorders the comma-separated "Recorded By" list in column G so that the
# "System" entries sort to the front of each cell (exact "System" first,
# then any other-case "system" variants), while any remaining names keep a
# stable alphabetical order. Mirrors a UI change where the session-analysis
# report started prioritising the automated "System" recorder in the
# display string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ", "
    if ($parts.Length -le 1) {
        continue
    }

    $sysExact = @()
    $sysCi = @()
    $others = @()
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $sysExact += $p
        } elseif ($p.ToLower().Equals("system")) {
            $sysCi += $p
        } else {
            $others += $p
        }
    }
    $others = $others | Sort-Object

    $newText = ($sysExact + $sysCi + $others) -join ", "

    if (-not $newText.Equals($text)) {
        $cell.Value = $newText
    }
}
